# update điểm của bảo duy

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (Gia Huy): mark as "Con chưa nộp bài" (reuses existing shared string) and
#     bump the row's custom height to match the wrapped note text.
$ws.Range("C4").Value = "Con chưa nộp bài"
$ws.Rows.Item(4).RowHeight = 42.75

# --- Row 8 (Bảo Duy): record the new grade + teacher note.
$ws.Range("B8").Value = 7.5

$note = "Duy làm bài tốt, con chú ý thêm các vấn đề sau:`n- Con cần làm thêm Bài 7.`n- Các bài toán tìm x cầm có kết luận ở cuối bài. "
$noteCell = $ws.Range("C8")
$noteCell.Value = $note

# Run 1 (chars 1-66): default formatting - left as-is.
# Run 2 (chars 67-76): " Bài 7.`n- " in bold.
$noteCell.Characters(67, 10).Font.Bold = $true
# Run 3 (chars 77-122): "Các bài toán tìm x cầm có kết luận ở cuối bài." regular weight.
$noteCell.Characters(77, 46).Font.Bold = $false
# Run 4 (char 123): trailing space in bold.
$noteCell.Characters(123, 1).Font.Bold = $true

$ws.Rows.Item(8).RowHeight = 54

# --- Selection, matching the author's cursor position when they saved.
[void]$ws.Range("G13").Select()
